# Fix presenter asterisks on final pres slides
$p = $ppt.ActivePresentation

# --- Slide 1: split the presenter-credits run so the asterisks move from
#     Travis Cox / Alex Hahn onto Cory Sabol / Josh Moore ---
$titleSlide = $p.Slides.Item(1)
$creditsShape = $titleSlide.Shapes.Item(2)
$creditsRange = $creditsShape.TextFrame.TextRange

$oldText = "Travis Cox*, Alex Hahn*, Cory Sabol, and Josh Moore"
$fullText = $creditsRange.Text
$startPos = $fullText.IndexOf($oldText) + 1

$firstPart = "Cory Sabol*, Josh Moore*, "
$secondPart = "Travis Cox, and Alex Hahn"

$firstRange = $creditsRange.Characters($startPos, $firstPart.Length)
$firstRange.Text = $firstPart

$secondRange = $creditsRange.Characters($startPos + $firstPart.Length, $secondPart.Length)
$secondRange.Text = $secondPart

# --- Slide 15: swap both tables over to the new table style ---
$newStyleId = "{8A2E84FA-CCA0-47B5-8583-5066F5180257}"
$tablesSlide = $p.Slides.Item(15)
for ($i = 1; $i -le $tablesSlide.Shapes.Count; $i++) {
    $shape = $tablesSlide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyleId)
    }
}
